# Apply the changes described by the diff:
# 1. Rename the sheet from "alpha3F-HW15.xpc" to "alpha3F"
# 2. Update a handful of row-15 values to their new (last-ULP-adjusted) doubles

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "alpha3F"

$ws.Range("D15").Value2 = 1.029101065489702
$ws.Range("I15").Value2 = 0.9698016854779219
$ws.Range("J15").Value2 = 1.029101065489702
$ws.Range("K15").Value2 = 1.004061216628706
$ws.Range("L15").Value2 = 0.9827844711854847
$ws.Range("M15").Value2 = 0.9832546231280049
